# The prior 'merge' (instead of 'join') step silently dropped any track that
# wasn't already present in every source, and the By_Artist_* views only ever
# credited the first artist of a multi-artist collab. This rewrites all nine
# report tabs with the complete, correctly-joined dataset (7 tracks / 9 artists
# / 4 labels), fixes the 'Loot At The Sky' -> 'Look At The Sky' typo, and splits
# 'Matt Nash, Lucas Marx' into its two artists for the By_Artist_* tabs.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: By_Track_YouTube ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = 'Porter Robinson'
$ws.Cells.Item(2, 2).Value = 'Look At The Sky'
$ws.Cells.Item(2, 3).Value = 'MOM+POP'
$ws.Cells.Item(2, 4).Value = 1378335
$ws.Cells.Item(3, 1).Value = 'Kaskade'
$ws.Cells.Item(3, 2).Value = 'Closer'
$ws.Cells.Item(3, 3).Value = 'Monstercat'
$ws.Cells.Item(3, 4).Value = 210865
$ws.Cells.Item(4, 1).Value = 'Rameses B'
$ws.Cells.Item(4, 2).Value = 'Samurai'
$ws.Cells.Item(4, 3).Value = 'Monstercat'
$ws.Cells.Item(4, 4).Value = 62992
$ws.Cells.Item(5, 1).Value = 'Matt Nash, Lucas Marx'
$ws.Cells.Item(5, 2).Value = 'Midnight'
$ws.Cells.Item(5, 3).Value = 'STMPD RCRDS'
$ws.Cells.Item(5, 4).Value = 41797
$ws.Cells.Item(6, 1).Value = 'Dirty Palm'
$ws.Cells.Item(6, 2).Value = 'Diamonds'
$ws.Cells.Item(6, 3).Value = 'NONE'
$ws.Cells.Item(6, 4).Value = 37749
$ws.Cells.Item(7, 1).Value = 'Body Ocean'
$ws.Cells.Item(7, 2).Value = 'Once The Music'
$ws.Cells.Item(7, 3).Value = 'STMPD RCRDS'
$ws.Cells.Item(7, 4).Value = 20803
$ws.Cells.Item(8, 1).Value = 'Lady Bee, Dame1'
$ws.Cells.Item(8, 2).Value = 'Soon Not Later'
$ws.Cells.Item(8, 3).Value = 'Mixmash Recorsds'
$ws.Cells.Item(8, 4).Value = 1958

# --- Sheet 2: By_Track_1001Tracklists ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 1).Value = 'Body Ocean'
$ws.Cells.Item(2, 2).Value = 'Once The Music'
$ws.Cells.Item(2, 3).Value = 'STMPD RCRDS'
$ws.Cells.Item(2, 4).Value = 64
$ws.Cells.Item(2, 5).Value = 66
$ws.Cells.Item(3, 1).Value = 'Matt Nash, Lucas Marx'
$ws.Cells.Item(3, 2).Value = 'Midnight'
$ws.Cells.Item(3, 3).Value = 'STMPD RCRDS'
$ws.Cells.Item(3, 4).Value = 48
$ws.Cells.Item(3, 5).Value = 53
$ws.Cells.Item(4, 1).Value = 'Kaskade'
$ws.Cells.Item(4, 2).Value = 'Closer'
$ws.Cells.Item(4, 3).Value = 'Monstercat'
$ws.Cells.Item(4, 4).Value = 8
$ws.Cells.Item(4, 5).Value = 11
$ws.Cells.Item(5, 1).Value = 'Porter Robinson'
$ws.Cells.Item(5, 2).Value = 'Look At The Sky'
$ws.Cells.Item(5, 3).Value = 'MOM+POP'
$ws.Cells.Item(5, 4).Value = 8
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(6, 1).Value = 'Lady Bee, Dame1'
$ws.Cells.Item(6, 2).Value = 'Soon Not Later'
$ws.Cells.Item(6, 3).Value = 'Mixmash Recorsds'
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(7, 1).Value = 'Dirty Palm'
$ws.Cells.Item(7, 2).Value = 'Diamonds'
$ws.Cells.Item(7, 3).Value = 'NONE'
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(8, 1).Value = 'Rameses B'
$ws.Cells.Item(8, 2).Value = 'Samurai'
$ws.Cells.Item(8, 3).Value = 'Monstercat'
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0

# --- Sheet 3: By_Track_Soundcloud ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 1).Value = 'Porter Robinson'
$ws.Cells.Item(2, 2).Value = 'Look At The Sky'
$ws.Cells.Item(2, 3).Value = 'MOM+POP'
$ws.Cells.Item(2, 4).Value = 79423
$ws.Cells.Item(3, 1).Value = 'Kaskade'
$ws.Cells.Item(3, 2).Value = 'Closer'
$ws.Cells.Item(3, 3).Value = 'Monstercat'
$ws.Cells.Item(3, 4).Value = 51258
$ws.Cells.Item(4, 1).Value = 'Rameses B'
$ws.Cells.Item(4, 2).Value = 'Samurai'
$ws.Cells.Item(4, 3).Value = 'Monstercat'
$ws.Cells.Item(4, 4).Value = 33843
$ws.Cells.Item(5, 1).Value = 'Dirty Palm'
$ws.Cells.Item(5, 2).Value = 'Diamonds'
$ws.Cells.Item(5, 3).Value = 'NONE'
$ws.Cells.Item(5, 4).Value = 9312
$ws.Cells.Item(6, 1).Value = 'Matt Nash, Lucas Marx'
$ws.Cells.Item(6, 2).Value = 'Midnight'
$ws.Cells.Item(6, 3).Value = 'STMPD RCRDS'
$ws.Cells.Item(6, 4).Value = 5967
$ws.Cells.Item(7, 1).Value = 'Body Ocean'
$ws.Cells.Item(7, 2).Value = 'Once The Music'
$ws.Cells.Item(7, 3).Value = 'STMPD RCRDS'
$ws.Cells.Item(7, 4).Value = 3477
$ws.Cells.Item(8, 1).Value = 'Lady Bee, Dame1'
$ws.Cells.Item(8, 2).Value = 'Soon Not Later'
$ws.Cells.Item(8, 3).Value = 'Mixmash Recorsds'
$ws.Cells.Item(8, 4).Value = 2924

# --- Sheet 4: By_Artist_YouTube ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 1).Value = 'Porter Robinson'
$ws.Cells.Item(2, 2).Value = 1378335
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 1).Value = 'Kaskade'
$ws.Cells.Item(3, 2).Value = 210865
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 1).Value = 'Rameses B'
$ws.Cells.Item(4, 2).Value = 62992
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(5, 1).Value = 'Lucas Marx'
$ws.Cells.Item(5, 2).Value = 41797
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(6, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(6, 1).Value = 'Matt Nash'
$ws.Cells.Item(6, 2).Value = 41797
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(7, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(7, 1).Value = 'Dirty Palm'
$ws.Cells.Item(7, 2).Value = 37749
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 1).Value = 'Body Ocean'
$ws.Cells.Item(8, 2).Value = 20803
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(9, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(9, 1).Value = 'Dame1'
$ws.Cells.Item(9, 2).Value = 1958
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(10, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10, 1).Value = 'Lady Bee'
$ws.Cells.Item(10, 2).Value = 1958

# --- Sheet 5: By_Artist_1001Tracklists ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 1).Value = 'Body Ocean'
$ws.Cells.Item(2, 2).Value = 64
$ws.Cells.Item(2, 3).Value = 66
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 1).Value = 'Lucas Marx'
$ws.Cells.Item(3, 2).Value = 48
$ws.Cells.Item(3, 3).Value = 53
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 1).Value = 'Matt Nash'
$ws.Cells.Item(4, 2).Value = 48
$ws.Cells.Item(4, 3).Value = 53
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(5, 1).Value = 'Kaskade'
$ws.Cells.Item(5, 2).Value = 8
$ws.Cells.Item(5, 3).Value = 11
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(6, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(6, 1).Value = 'Porter Robinson'
$ws.Cells.Item(6, 2).Value = 8
$ws.Cells.Item(6, 3).Value = 9
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(7, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(7, 1).Value = 'Dame1'
$ws.Cells.Item(7, 2).Value = 3
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 1).Value = 'Lady Bee'
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(8, 3).Value = 3
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(9, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(9, 1).Value = 'Dirty Palm'
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(10, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10, 1).Value = 'Rameses B'
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 0

# --- Sheet 6: By_Artist_Soundcloud ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 1).Value = 'Porter Robinson'
$ws.Cells.Item(2, 2).Value = 79423
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 1).Value = 'Kaskade'
$ws.Cells.Item(3, 2).Value = 51258
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 1).Value = 'Rameses B'
$ws.Cells.Item(4, 2).Value = 33843
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(5, 1).Value = 'Dirty Palm'
$ws.Cells.Item(5, 2).Value = 9312
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(6, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(6, 1).Value = 'Lucas Marx'
$ws.Cells.Item(6, 2).Value = 5967
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(7, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(7, 1).Value = 'Matt Nash'
$ws.Cells.Item(7, 2).Value = 5967
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(8, 1).Value = 'Body Ocean'
$ws.Cells.Item(8, 2).Value = 3477
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(9, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(9, 1).Value = 'Dame1'
$ws.Cells.Item(9, 2).Value = 2924
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(10, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10, 1).Value = 'Lady Bee'
$ws.Cells.Item(10, 2).Value = 2924

# --- Sheet 7: By_Label_YouTube ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(2, 1).Value = 'MOM+POP'
$ws.Cells.Item(2, 2).Value = 1378335
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 1).Value = 'Monstercat'
$ws.Cells.Item(3, 2).Value = 273857
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 1).Value = 'STMPD RCRDS'
$ws.Cells.Item(4, 2).Value = 62600
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(5, 1).Value = 'Mixmash Recorsds'
$ws.Cells.Item(5, 2).Value = 1958

# --- Sheet 8: By_Label_1001Tracklists ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2, 1).Value = 'STMPD RCRDS'
$ws.Cells.Item(2, 2).Value = 112
$ws.Cells.Item(2, 3).Value = 119
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 1).Value = 'Monstercat'
$ws.Cells.Item(3, 2).Value = 8
$ws.Cells.Item(3, 3).Value = 11
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 1).Value = 'MOM+POP'
$ws.Cells.Item(4, 2).Value = 8
$ws.Cells.Item(4, 3).Value = 9
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(5, 1).Value = 'Mixmash Recorsds'
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = 3

# --- Sheet 9: By_Label_Soundcloud ---
$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 1).Value = 'Monstercat'
$ws.Cells.Item(2, 2).Value = 85101
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(3, 1).Value = 'MOM+POP'
$ws.Cells.Item(3, 2).Value = 79423
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 1).Value = 'STMPD RCRDS'
$ws.Cells.Item(4, 2).Value = 9444
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(5, 1).Value = 'Mixmash Recorsds'
$ws.Cells.Item(5, 2).Value = 2924

$excel.CutCopyMode = 0
